$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 1.03
$ws.Range("Q2").Value = 1.61
$ws.Range("R2").Value = 1.56
$ws.Range("T2").Value = 1.56
$ws.Range("X2").Value = 27
$ws.Range("AB2").Value = 19.5
$ws.Range("AK2").Value = 34

# Row 5
$ws.Range("I5").Value = 1.61

# Row 8
$ws.Range("H8").Value = 7.4
$ws.Range("I8").Value = 7.6
$ws.Range("J8").Value = 5.1
$ws.Range("K8").Value = 5.2
$ws.Range("N8").Value = 5.8
$ws.Range("P8").Value = 2.62
$ws.Range("Q8").Value = 1.6
$ws.Range("R8").Value = 1.64
$ws.Range("S8").Value = 2.46
$ws.Range("AF8").Value = 9.800000000000001
$ws.Range("AL8").Value = 28
$ws.Range("AN8").Value = 5.6

# Row 10
$ws.Range("H10").Value = 3.85
$ws.Range("K10").Value = 3.65
$ws.Range("L10").Value = 1.37
$ws.Range("Q10").Value = 2.04
$ws.Range("W10").Value = 1.84
$ws.Range("AF10").Value = 15.5
$ws.Range("AN10").Value = 23

# Row 11
$ws.Range("F11").Value = 1.85
$ws.Range("H11").Value = 3.65
$ws.Range("I11").Value = 6
$ws.Range("K11").Value = 5.4
$ws.Range("Q11").Value = 1.96
